# Update report for vref=250mV
# Updates offset_val (R), offset (S), and comp_resolution (T) columns
# for the rows whose measured values changed with the new reference voltage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("R2").Value = 250.514
$ws.Range("S2").Value = -99.48650000000001
$ws.Range("T2").Value = -0.513522

# Row 5
$ws.Range("R5").Value = 250.743
$ws.Range("S5").Value = -64.2572
$ws.Range("T5").Value = -0.742795

# Row 6
$ws.Range("R6").Value = 267.523
$ws.Range("S6").Value = -47.4774
$ws.Range("T6").Value = -17.5226

# Row 105
$ws.Range("R105").Value = 250.514
$ws.Range("S105").Value = -99.48650000000001
$ws.Range("T105").Value = -17.5226

# Row 106
$ws.Range("R106").Value = 267.523
$ws.Range("S106").Value = -47.4774
$ws.Range("T106").Value = -0.513522

# Row 111
$ws.Range("R111").Value = 256.26
$ws.Range("S111").Value = -70.40703000000001
$ws.Range("T111").Value = -6.259639

# Row 112
$ws.Range("R112").Value = 7.964692
$ws.Range("S112").Value = 21.673362
$ws.Range("T112").Value = 7.964666

# Row 113
$ws.Range("R113").Value = 3.108051
$ws.Range("S113").Value = 30.782952
$ws.Range("T113").Value = 127.238424
